# Update "paises.xlsx" ("Pais" sheet) with the latest COVID-19 snapshot.
# - Refresh the "last updated" timestamp.
# - Refresh raw case counters for a handful of countries.
# - Panama's case count overtook Dinamarca / Serbia / Corea del Sur, so the
#   ranked table reshuffles: Panama's row moves up (row 49) with fresh
#   data, and the three countries it passed shift down one row each,
#   keeping their previous totals.
# - Groenlandia / Islas Turcas y Caicos swap ranking order as well (rows
#   207-208), each keeping its own (unchanged) totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 02:05"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Range("B4").Value = 1725400
$ws.Range("C4").Value = 19174
$ws.Range("D4").Value = 478237
$ws.Range("E4").Value = 1146618
$ws.Range("G4").Value = 740
$ws.Range("H4").Value = 100545

# --- Row 5: Brasil ------------------------------------------------------
$ws.Range("B5").Value = 392360
$ws.Range("C5").Value = 15691
$ws.Range("E5").Value = 209218
$ws.Range("G5").Value = 1027
$ws.Range("H5").Value = 24549

# --- Row 16: Turquia ------------------------------------------------
$ws.Range("B16").Value = 86647
$ws.Range("C16").Value = 936
$ws.Range("D16").Value = 45339
$ws.Range("E16").Value = 34669

# --- Rows 49-52: Panama overtakes Dinamarca / Serbia / Corea del Sur ---
# Row 49 becomes Panama with its new totals.
$ws.Range("A49").Value = "Panama"
$ws.Range("B49").Value = 11447
$ws.Range("C49").Value = 264
$ws.Range("D49").Value = 6379
$ws.Range("E49").Value = 4755
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 313

# Row 50 becomes Dinamarca, keeping its previous totals.
$ws.Range("A50").Value = "Dinamarca"
$ws.Range("B50").Value = 11428
$ws.Range("C50").Value = 41
$ws.Range("D50").Value = 10044
$ws.Range("E50").Value = 821
$ws.Range("H50").Value = 563

# Row 51 becomes Serbia, keeping its previous totals.
$ws.Range("A51").Value = "Serbia"
$ws.Range("B51").Value = 11227
$ws.Range("C51").Value = 34
$ws.Range("D51").Value = 6067
$ws.Range("E51").Value = 4921
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 239

# Row 52 becomes Corea del Sur, keeping its previous totals.
$ws.Range("A52").Value = "Corea del Sur"
$ws.Range("B52").Value = 11225
$ws.Range("C52").Value = 19
$ws.Range("D52").Value = 10275
$ws.Range("E52").Value = 681
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 269

# --- Row 65: Azerbaiyan -------------------------------------------------
$ws.Range("B65").Value = 7117
$ws.Range("C65").Value = 309
$ws.Range("D65").Value = 2317
$ws.Range("E65").Value = 4766
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 34

# --- Row 165 ------------------------------------------------------------
$ws.Range("D165").Value = 91
$ws.Range("E165").Value = 39

# --- Row 166 --------------------------------------------------------
$ws.Range("B166").Value = 137
$ws.Range("C166").Value = 3
$ws.Range("D166").Value = 63
$ws.Range("E166").Value = 73

# --- Rows 207-208: Groenlandia / Islas Turcas y Caicos swap order -------
# Row 207 becomes Groenlandia, keeping its previous totals.
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("D207").Value = 11
$ws.Range("H207").Value = 0

# Row 208 becomes Islas Turcas y Caicos, keeping its previous totals.
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 10
$ws.Range("H208").Value = 1
